$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value2 = 154.90909
$ws.Range("I11").Value2 = 154.90909
$ws.Range("K11").Value2 = 154.90909
$ws.Range("M11").Value2 = -14.90908999999999
# Row 19
$ws.Range("H19").Value2 = 1008.2
$ws.Range("J19").Value2 = 855.6667
$ws.Range("L19").Value2 = 855.6667
$ws.Range("N19").Value2 = -1205.6667
# Row 39
$ws.Range("H39").Value2 = 328.6
$ws.Range("I39").Value2 = 328.6
$ws.Range("K39").Value2 = 985.8000000000001
$ws.Range("M39").Value2 = -689.8000000000001
# Row 64
$ws.Range("H64").Value2 = 3994
$ws.Range("I64").Value2 = 3661
$ws.Range("K64").Value2 = 3661
$ws.Range("M64").Value2 = -3413
# Row 67
$ws.Range("H67").Value2 = 3994
$ws.Range("I67").Value2 = 3661
$ws.Range("K67").Value2 = 3661
$ws.Range("M67").Value2 = -2803
# Row 74
$ws.Range("H74").Value2 = 4767.8335
$ws.Range("I74").Value2 = 4767.8335
$ws.Range("K74").Value2 = 4767.8335
$ws.Range("M74").Value2 = -3831.8335
# Row 77
$ws.Range("H77").Value2 = 4767.8335
$ws.Range("I77").Value2 = 4767.8335
$ws.Range("K77").Value2 = 23839.1675
$ws.Range("M77").Value2 = -19159.1675
# Row 125
$ws.Range("H125").Value2 = 800.65
# Row 141
$ws.Range("H141").Value2 = 2940.5
$ws.Range("I141").Value2 = 2940.5
$ws.Range("K141").Value2 = 8821.5
$ws.Range("M141").Value2 = -3641.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 76
$ws.Range("H76").Value2 = 52037.4
$ws.Range("J76").Value2 = 52037.4
$ws.Range("L76").Value2 = 52037.4
$ws.Range("N76").Value2 = -52713.4
# Row 79
$ws.Range("H79").Value2 = 52037.4
$ws.Range("J79").Value2 = 52037.4
$ws.Range("L79").Value2 = 52037.4
$ws.Range("N79").Value2 = -54377.4
# Row 102
$ws.Range("H102").Value2 = 3195.6365
$ws.Range("I102").Value2 = 3128
$ws.Range("K102").Value2 = 3128
$ws.Range("M102").Value2 = -1506
# Row 132
$ws.Range("H132").Value2 = 5624.5
$ws.Range("I132").Value2 = 4463.4062
$ws.Range("K132").Value2 = 13390.2186
$ws.Range("M132").Value2 = -10860.2186

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 44
$ws.Range("H44").Value2 = 0
$ws.Range("J44").Value2 = 0
$ws.Range("L44").Value2 = 0
$ws.Range("N44").ClearContents()
# Row 107
$ws.Range("H107").Value2 = 3404.0833
$ws.Range("I107").Value2 = 2841.5
$ws.Range("K107").Value2 = 2841.5
$ws.Range("M107").Value2 = -921.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value2 = 499.57144
$ws.Range("I16").Value2 = 422.6154
$ws.Range("K16").Value2 = 422.6154
$ws.Range("M16").Value2 = -135.6154
# Row 32
$ws.Range("H32").Value2 = 25000
$ws.Range("I32").Value2 = 0
$ws.Range("J32").Value2 = 25000
$ws.Range("K32").Value2 = 0
$ws.Range("L32").Value2 = 25000
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value2 = -25632
# Row 35
$ws.Range("H35").Value2 = 1553.4286
$ws.Range("I35").Value2 = 1479
$ws.Range("J35").Value2 = 2000
$ws.Range("K35").Value2 = 1479
$ws.Range("L35").Value2 = 2000
$ws.Range("M35").Value2 = -1185
$ws.Range("N35").Value2 = -2588
# Row 39
$ws.Range("H39").Value2 = 4999
$ws.Range("I39").Value2 = 4999
$ws.Range("J39").Value2 = 0
$ws.Range("K39").Value2 = 4999
$ws.Range("L39").Value2 = 0
$ws.Range("M39").Value2 = -4608
$ws.Range("N39").ClearContents()
# Row 49
$ws.Range("H49").Value2 = 4999
$ws.Range("I49").Value2 = 4999
$ws.Range("J49").Value2 = 0
$ws.Range("K49").Value2 = 4999
$ws.Range("L49").Value2 = 0
$ws.Range("M49").Value2 = -4817
$ws.Range("N49").ClearContents()
# Row 51
$ws.Range("H51").Value2 = 49999
$ws.Range("I51").Value2 = 39999
$ws.Range("J51").Value2 = 69999
$ws.Range("K51").Value2 = 39999
$ws.Range("L51").Value2 = 69999
$ws.Range("M51").Value2 = -39263
$ws.Range("N51").Value2 = -71471
# Row 61
$ws.Range("H61").Value2 = 49999
$ws.Range("I61").Value2 = 39999
$ws.Range("J61").Value2 = 69999
$ws.Range("K61").Value2 = 39999
$ws.Range("L61").Value2 = 69999
$ws.Range("M61").Value2 = -39651
$ws.Range("N61").Value2 = -70695
# Row 62
$ws.Range("H62").Value2 = 4000
$ws.Range("I62").Value2 = 2500
$ws.Range("J62").Value2 = 4750
$ws.Range("K62").Value2 = 2500
$ws.Range("L62").Value2 = 4750
$ws.Range("M62").Value2 = -1876
$ws.Range("N62").Value2 = -5998
# Row 65
$ws.Range("H65").Value2 = 4000
$ws.Range("I65").Value2 = 2500
$ws.Range("J65").Value2 = 4750
$ws.Range("K65").Value2 = 12500
$ws.Range("L65").Value2 = 23750
$ws.Range("M65").Value2 = -9380
$ws.Range("N65").Value2 = -29990
# Row 94
$ws.Range("J94").Value2 = 1200
$ws.Range("L94").Value2 = 1200
$ws.Range("N94").Value2 = -2102
# Row 105
$ws.Range("H105").Value2 = 1795
$ws.Range("I105").Value2 = 1488
$ws.Range("K105").Value2 = 1488
$ws.Range("M105").Value2 = 259
# Row 113
$ws.Range("H113").Value2 = 499.57144
$ws.Range("I113").Value2 = 422.6154
$ws.Range("K113").Value2 = 422.6154
$ws.Range("M113").Value2 = 1747.3846

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value2 = 1139.6
$ws.Range("I113").Value2 = 1152.5652
$ws.Range("J113").Value2 = 990.5
$ws.Range("K113").Value2 = 3457.6956
$ws.Range("L113").Value2 = 2971.5
$ws.Range("M113").Value2 = -1287.6956
$ws.Range("N113").Value2 = -7311.5
# Row 122
$ws.Range("H122").Value2 = 713.1429000000001
$ws.Range("J122").Value2 = 759
$ws.Range("L122").Value2 = 6831
$ws.Range("N122").Value2 = -11731
# Row 124
$ws.Range("H124").Value2 = 199
$ws.Range("I124").Value2 = 199
$ws.Range("J124").Value2 = 0
$ws.Range("K124").Value2 = 597
$ws.Range("L124").Value2 = 0
$ws.Range("M124").Value2 = 4313
$ws.Range("N124").ClearContents()
# Row 131
$ws.Range("H131").Value2 = 2324.6667
$ws.Range("I131").Value2 = 2020
$ws.Range("J131").Value2 = 2426.2222
$ws.Range("K131").Value2 = 6060
$ws.Range("L131").Value2 = 7278.6666
$ws.Range("M131").Value2 = -1020
$ws.Range("N131").Value2 = -17358.6666
# Row 132
$ws.Range("H132").Value2 = 1000
$ws.Range("I132").Value2 = 1000
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 9000
$ws.Range("L132").Value2 = 0
$ws.Range("M132").Value2 = -6470
$ws.Range("N132").ClearContents()
# Row 139
$ws.Range("H139").Value2 = 6674.75
$ws.Range("I139").Value2 = 2457.5
$ws.Range("J139").Value2 = 8080.5
$ws.Range("K139").Value2 = 7372.5
$ws.Range("L139").Value2 = 24241.5
$ws.Range("M139").Value2 = -2232.5
$ws.Range("N139").Value2 = -34521.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value2 = 0
$ws.Range("J70").Value2 = 0
$ws.Range("L70").Value2 = 0
$ws.Range("N70").ClearContents()
# Row 73
$ws.Range("H73").Value2 = 0
$ws.Range("J73").Value2 = 0
$ws.Range("L73").Value2 = 0
$ws.Range("N73").ClearContents()
# Row 102
$ws.Range("H102").Value2 = 1557.0358
$ws.Range("I102").Value2 = 1278.6086
$ws.Range("K102").Value2 = 1278.6086
$ws.Range("M102").Value2 = 343.3914
# Row 126
$ws.Range("H126").Value2 = 4057
$ws.Range("I126").Value2 = 1100
$ws.Range("K126").Value2 = 3300
$ws.Range("M126").Value2 = -830

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 32
$ws.Range("H32").Value2 = 6055.75
$ws.Range("I32").Value2 = 2361.5
$ws.Range("J32").Value2 = 9750
$ws.Range("K32").Value2 = 2361.5
$ws.Range("L32").Value2 = 9750
$ws.Range("M32").Value2 = -2044.5
$ws.Range("N32").Value2 = -10384
# Row 68
$ws.Range("H68").Value2 = 3237.25
$ws.Range("I68").Value2 = 3316.3333
$ws.Range("K68").Value2 = 3316.3333
$ws.Range("M68").Value2 = -2567.3333
# Row 71
$ws.Range("H71").Value2 = 3237.25
$ws.Range("I71").Value2 = 3316.3333
$ws.Range("K71").Value2 = 16581.6665
$ws.Range("M71").Value2 = -12837.6665
# Row 93
$ws.Range("H93").Value2 = 0
$ws.Range("J93").Value2 = 0
$ws.Range("L93").Value2 = 0
$ws.Range("N93").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value2 = 15502.259
$ws.Range("I136").Value2 = 14778
$ws.Range("K136").Value2 = 44334
$ws.Range("M136").Value2 = -41784
